$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths (target OOXML width 15.42578125 for both columns;
# Excel quantizes ColumnWidth to whole pixels, so 15.5 is the nearest representable width)
$ws.Columns.Item(1).ColumnWidth = 14.7109375
$ws.Columns.Item(2).ColumnWidth = 14.7109375

$ws.Range("A1").Value = -0.38267251156081272
$ws.Range("B1").Value = 0.38181338606057125
$ws.Range("A2").Value = -0.33570347490350194
$ws.Range("B2").Value = 0.33214555221213615
$ws.Range("A3").Value = -0.17429041547444868
$ws.Range("B3").Value = 0.17357622783611149
$ws.Range("A4").Value = -0.16157622801121363
$ws.Range("B4").Value = 0.1609350299963026
$ws.Range("A5").Value = -0.15493503063787362
$ws.Range("B5").Value = 0.15365376561521149
$ws.Range("A6").Value = -0.055280449236837281
$ws.Range("B6").Value = 0.055233441416508544
$ws.Range("A7").Value = -0.035233442191284325
$ws.Range("B7").Value = 0.035162218168169446
$ws.Range("A8").Value = -0.015162218946647599
$ws.Range("B8").Value = 0.015143370119504951
$ws.Range("A9").Value = -0.0091433707893058269
$ws.Range("B9").Value = 0.0091361839629859887
$ws.Range("A10").Value = -0.0031361846341582122
$ws.Range("B10").Value = 0.0031378898906737618
$ws.Range("A11").Value = 0.0013621094500493314
$ws.Range("B11").Value = -0.0013612077917244392
$ws.Range("A12").Value = 0.007361207120549107
$ws.Range("B12").Value = -0.0073675050806336628
$ws.Range("A13").Value = 0.013367504411212927
$ws.Range("B13").Value = -0.013375120472231217
$ws.Range("A14").Value = 0.02451394357644876
$ws.Range("B14").Value = -0.024537070742218781
$ws.Range("A15").Value = 0.030537070076142925
$ws.Range("B15").Value = -0.030583410637460062
$ws.Range("A16").Value = -0.015026457675629601
$ws.Range("B16").Value = 0.015003896376523329
$ws.Range("A17").Value = -0.0090038970421044695
$ws.Range("B17").Value = 0.0089999993097231723
$ws.Range("A18").Value = -0.036109911785644044
$ws.Range("B18").Value = 0.036096109618934946
$ws.Range("A19").Value = -0.027096110237945226
$ws.Range("B19").Value = 0.027013014072142294
$ws.Range("A20").Value = -0.018013014697165985
$ws.Range("B20").Value = 0.018004086428785371
$ws.Range("A21").Value = -0.0090040870548335761
$ws.Range("B21").Value = 0.0089999993731506578
$ws.Range("A22").Value = -0.093951177780187578
$ws.Range("B22").Value = 0.093637264880840121
$ws.Range("A23").Value = -0.084637265532695238
$ws.Range("B23").Value = 0.084127223164658815
$ws.Range("A24").Value = -0.042127224089963278
$ws.Range("B24").Value = 0.041999999068885074
$ws.Range("A25").Value = -0.092373317040465253
$ws.Range("B25").Value = 0.092253205258678861
$ws.Range("A26").Value = -0.085222882449691895
$ws.Range("B26").Value = 0.085053589654364714
$ws.Range("A27").Value = -0.079053590322519796
$ws.Range("B27").Value = 0.078492731005013194
$ws.Range("A28").Value = -0.072492731687163747
$ws.Range("B28").Value = 0.072125231875200591
$ws.Range("A29").Value = -0.060125232614856472
$ws.Range("B29").Value = 0.059967173504889004
$ws.Range("A30").Value = -0.03996717431304786
$ws.Range("B30").Value = 0.039804729553190032
$ws.Range("A31").Value = -0.024804730328533253
$ws.Range("B31").Value = 0.02476556800306895
$ws.Range("A32").Value = -0.0060004922835927843
$ws.Range("B32").Value = 0.0059999992947021852
